# Add "compress ratio (iteration = 1000, clusters = 1000)" table to the
# k-means sheet, below the existing "execution time (iteration = 1000,
# clusters = 1000)" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("k-means")

# Section title
$ws.Range("A18").Value = "compress ratio (iteration = 1000, clusters = 1000)"

# absErrBound header row
$ws.Range("A19").Value = "absErrBound"
$ws.Range("B19").Value = 0.000001
$ws.Range("C19").Value = 0.00001
$ws.Range("D19").Value = 0.0001
$ws.Range("E19").Value = 0.001
$ws.Range("F19").Value = 0.01
$ws.Range("G19").Value = 0.1

# compress ratio values
$ws.Range("B20").Formula = "= 1/0.979871"
$ws.Range("C20").Value = 1.019641
$ws.Range("D20").Value = 1.036257
$ws.Range("E20").Value = 1.059337
$ws.Range("F20").Value = 1.156326
$ws.Range("G20").Value = 1.359143

# Page setup for the sheet (new print area settings introduced with this edit)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

[void]$ws.Activate()
[void]$ws.Range("A18").Select()
